# Refresh the cryptocurrency price/volume table (cols B-E, rows 2-51) to match
# the latest scrape. Price cells in column D are numeric-looking text (e.g.
# "51.097.23", "0.110") so a leading apostrophe is used to force Excel to
# store them as literal text instead of re-parsing/re-formatting them as
# numbers (which would silently drop meaningful trailing/leading zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.097.23'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '''2.960.43'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''379.73'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').Value = '''102.27'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.590'
$ws.Range('E9').Value = '  +1.54%  '
$ws.Range('D10').Value = '''36.56'
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('E12').Value = '  +2.15%  '
$ws.Range('D13').Value = '''3.426.73'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '''7.80'
$ws.Range('E14').Value = '  +6.46%  '
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '''11.95'
$ws.Range('E16').Value = '  +67.08%  '
$ws.Range('D17').Value = '''2.968.11'
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').Value = '''0.999'
$ws.Range('E18').Value = '  +2.61%  '
$ws.Range('D19').Value = '''51.170.01'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('E21').Value = '  -0.85%  '
$ws.Range('D22').Value = '''0.0₃0960'
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('D23').Value = '''70.13'
$ws.Range('E23').Value = '  +2.73%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '''268.01'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '''3.27'
$ws.Range('E25').Value = '  +13.56%  '
$ws.Range('D26').Value = '''7.90'
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('D27').Value = '''7.19'
$ws.Range('E27').Value = '  -7.87%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').Value = '''25.89'
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('D31').Value = '''0.110'
$ws.Range('E31').Value = '  -1.75%  '
$ws.Range('E32').Value = '  +6.22%  '
$ws.Range('D33').Value = '''34.46'
$ws.Range('E33').Value = '  +2.91%  '
$ws.Range('D34').Value = '''51.05'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('E35').Value = '  +2.62%  '
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  +9.33%  '
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '''1.83'
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '''16.58'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '''125.23'
$ws.Range('E42').Value = '  +4.06%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''2.50'
$ws.Range('E43').Value = '  -1.97%  '
$ws.Range('D44').Value = '''21.61'
$ws.Range('E44').Value = '  +3.27%  '
$ws.Range('E45').Value = '  +9.51%  '
$ws.Range('E46').Value = '  +2.98%  '
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('D48').Value = '''2.052.44'
$ws.Range('E48').Value = '  +4.15%  '
$ws.Range('D49').Value = '''0.270'
$ws.Range('E49').Value = '  -6.19%  '
$ws.Range('E50').Value = '  -6.67%  '
$ws.Range('E51').Value = '  +7.71%  '
